$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("D5").Value = "칼만 필터"

# Row 9
$ws.Range("D9").Value = "일반 MBA/DBA 석사/박사 과정 커리큘럼"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/general-mba-dba-course-works/#utm_source=rss&utm_medium=rss&utm_campaign=general-mba-dba-course-works"

# Row 16
$ws.Range("D16").Value = "SmoothGrad : removing noise by adding noise 내용 정리 [XAI-5]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/146"

# Row 23
$ws.Range("D23").Value = "Keras에서 입력영상의 컬러채널 또는 다수의 깊이(depth channel)채널과 필터의 콘볼루션 영상 후 output이 어떻게 계산될까?"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2787"

# Row 28
$ws.Range("D28").Value = "PointCloud :: Range image란?"
$ws.Range("E28").Value = "https://ropiens.tistory.com/108"

# Row 29
$ws.Range("D29").Value = "[만화] 인턴일기 1~7"
$ws.Range("E29").Value = "https://blog.promedius.ai/intern-life-1/"

# Row 32
$ws.Range("D32").Value = "데커레이터 다수 지정할때 실행 순서?!"
$ws.Range("E32").Value = "https://dodonam.tistory.com/315"

# Row 37
$ws.Range("D37").Value = "[Paper Review] Deep GNNs"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1457&mod=document&pageid=1"

# Row 51
$ws.Range("D51").Value = "윈도우 실행창에서 명령어로 시스템 속성창 열기, sysdm.cpl"
$ws.Range("E51").Value = "https://bskyvision.com/1172"
